$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($row, $value) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-PriceText 2 "242.64"
Set-PriceText 3 "23.11"
Set-PriceText 4 "5.421"
Set-PriceText 5 "0.05895"
Set-PriceText 6 "3.433"
Set-PriceText 7 "6.532"
Set-PriceText 8 "0.8088"
Set-PriceText 9 "0.9311"
Set-PriceText 10 "0.1419"
Set-PriceText 11 "0.07366"
Set-PriceText 12 "0.03306"
Set-PriceText 13 "0.03068"
Set-PriceText 14 "0.09350"
Set-PriceText 15 "3.847"
Set-PriceText 16 "0.001573"
Set-PriceText 17 "0.04691"
Set-PriceText 18 "0.0005942"
Set-PriceText 19 "0.005981"
Set-PriceText 20 "0.001258"
Set-PriceText 21 "0.004901"
Set-PriceText 22 "0.00006800"
Set-PriceText 23 "3.564"
Set-PriceText 24 "2.143"
Set-PriceText 26 "0.1294"
Set-PriceText 40 "0.03974"
Set-PriceText 41 "0.006182"
Set-PriceText 42 "0.1072"
Set-PriceText 43 "0.003000"
Set-PriceText 44 "0.008741"
Set-PriceText 45 "0.00005330"
Set-PriceText 48 "0.002335"
Set-PriceText 49 "0.00002100"
Set-PriceText 50 "0.0002000"
